# Atualização de bases das ligas, do dia: 28-05-2024 às 07:50
#
# The source data rows got re-ordered / re-matched against their betting
# odds lines. The fix re-shuffles the content of a handful of data rows
# (columns B..AD) while leaving the running index in column A untouched:
#   - rows 59 / 60          -> swapped
#   - rows 108 / 109 / 110  -> rotated (108<-109, 109<-110, 110<-108)
#   - rows 226 / 227        -> swapped
#   - rows 232 / 233        -> swapped
#   - rows 236 / 237        -> swapped

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B..AD hold the data that needs to move; column A (the running
# index) must stay exactly where it is.
$firstCol = 2   # B
$lastCol  = 30  # AD

function Get-RowValues($row) {
    $vals = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $vals += $ws.Cells.Item($row, $c).Value()
    }
    return $vals
}

function Set-RowValues($row, $vals) {
    $i = 0
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($row, $c).Value = $vals[$i]
        $i++
    }
}

# --- swap helper -----------------------------------------------------
function Swap-Rows($rowA, $rowB) {
    $dataA = Get-RowValues $rowA
    $dataB = Get-RowValues $rowB
    Set-RowValues $rowA $dataB
    Set-RowValues $rowB $dataA
}

# --- 3-way rotation helper (target gets data from "source") ----------
function Rotate-Rows($r1, $r2, $r3) {
    # r1 <- r2, r2 <- r3, r3 <- r1 (original)
    $d1 = Get-RowValues $r1
    $d2 = Get-RowValues $r2
    $d3 = Get-RowValues $r3
    Set-RowValues $r1 $d2
    Set-RowValues $r2 $d3
    Set-RowValues $r3 $d1
}

Swap-Rows 59 60
Rotate-Rows 108 109 110
Swap-Rows 226 227
Swap-Rows 232 233
Swap-Rows 236 237
